$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.731.94"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.906.58"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'312.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'0.5209"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.19%  "
$ws.Range("D8").Value = "'0.3785"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "'0.07251"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").Value = "'21.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.16%  "
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").Value = "'0.07660"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "1.883.50"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").Value = "'5.447"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "'92.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "'0.000008712"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "27.775.07"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "'14.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").Value = "'5.146"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").Value = "2.141.03"
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("D23").Value = "'10.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("D24").Value = "'6.639"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "'153.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "'1.870"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").Value = "'2.169"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("D28").Value = "'18.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").Value = "'114.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").Value = "'4.860"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").Value = "'0.09093"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("D32").Value = "'3.189"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").Value = "'4.842"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.69%  "
$ws.Range("D34").Value = "'1.230"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D35").Value = "'0.7795"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.11%  "
$ws.Range("E36").Value = "  +2.48%  "
$ws.Range("D37").Value = "'2.602"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.04%  "
$ws.Range("D38").Value = "'3.074"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.06%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.5587"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.52%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.093"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").Value = "'0.05284"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").Value = "'6.714"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.54%  "
$ws.Range("D43").Value = "'114.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.12%  "
$ws.Range("D44").Value = "'8.542"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").Value = "'0.4815"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("D48").Value = "'1.000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("D50").Value = "'66.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").Value = "'0.05989"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.11%  "
